$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 11
$ws.Range("I6").Value = 11
$ws.Range("K6").Value = 33
$ws.Range("M6").Value = 79
$ws.Range("H9").Value = 838.2857
$ws.Range("I9").Value = 946.3333
$ws.Range("K9").Value = 946.3333
$ws.Range("M9").Value = -777.3333
$ws.Range("H17").Value = 2996.6667
$ws.Range("J17").Value = 2996.6667
$ws.Range("L17").Value = 8990.000100000001
$ws.Range("N17").Value = -9326.000100000001
$ws.Range("H51").Value = 10061.6
$ws.Range("I51").Value = 8792.6
$ws.Range("K51").Value = 8792.6
$ws.Range("M51").Value = -8308.6
$ws.Range("H54").Value = 28000
$ws.Range("I54").Value = 28000
$ws.Range("K54").Value = 28000
$ws.Range("M54").Value = -27514
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H113").Value = 3134.4
$ws.Range("I113").Value = 3094.3076
$ws.Range("K113").Value = 3094.3076
$ws.Range("M113").Value = 159.6923999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1781
$ws.Range("I45").Value = 1468.8334
$ws.Range("K45").Value = 1468.8334
$ws.Range("M45").Value = -1091.8334
$ws.Range("H61").Value = 1594.6428
$ws.Range("I61").Value = 1409.6923
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 1409.6923
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -1197.6923
$ws.Range("N61").Value = -4423
$ws.Range("H97").Value = 2427.6
$ws.Range("I97").Value = 1301.0769
$ws.Range("K97").Value = 1301.0769
$ws.Range("M97").Value = -805.0769
$ws.Range("H136").Value = 1594.6428
$ws.Range("I136").Value = 1409.6923
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 4229.0769
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -1679.0769
$ws.Range("N136").Value = -17097

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H94").Value = 4526.9287
$ws.Range("I94").Value = 4184.1665
$ws.Range("J94").Value = 4784
$ws.Range("K94").Value = 4184.1665
$ws.Range("L94").Value = 4784
$ws.Range("M94").Value = -3733.1665
$ws.Range("N94").Value = -5686
$ws.Range("H105").Value = 3148.625
$ws.Range("I105").Value = 2538
$ws.Range("J105").Value = 4166.3335
$ws.Range("K105").Value = 2538
$ws.Range("L105").Value = 4166.3335
$ws.Range("M105").Value = -791
$ws.Range("N105").Value = -7660.3335
$ws.Range("H123").Value = 49990
$ws.Range("J123").Value = 49990
$ws.Range("L123").Value = 49990
$ws.Range("N123").Value = -59790

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 17821.25
$ws.Range("J96").Value = 17821.25
$ws.Range("L96").Value = 17821.25
$ws.Range("N96").Value = -23313.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2691.125
$ws.Range("I20").Value = 4507.5
$ws.Range("J20").Value = 874.75
$ws.Range("K20").Value = 13522.5
$ws.Range("L20").Value = 2624.25
$ws.Range("M20").Value = -13295.5
$ws.Range("N20").Value = -3078.25
$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224
$ws.Range("H68").Value = 297.5
$ws.Range("I68").Value = 297.5
$ws.Range("K68").Value = 892.5
$ws.Range("M68").Value = -81.5
$ws.Range("H71").Value = 297.5
$ws.Range("I71").Value = 297.5
$ws.Range("K71").Value = 2677.5
$ws.Range("M71").Value = 1378.5
$ws.Range("H131").Value = 1819.2632
$ws.Range("J131").Value = 1783.5625
$ws.Range("L131").Value = 5350.6875
$ws.Range("N131").Value = -15430.6875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 23668.5
$ws.Range("J95").Value = 23668.5
$ws.Range("L95").Value = 23668.5
$ws.Range("N95").Value = -29160.5
$ws.Range("H97").Value = 998.2222
$ws.Range("I97").Value = 248
$ws.Range("J97").Value = 7000
$ws.Range("K97").Value = 248
$ws.Range("L97").Value = 7000
$ws.Range("M97").Value = 248
$ws.Range("N97").Value = -7992
$ws.Range("H107").Value = 1674.1538
$ws.Range("J107").Value = 4741.25
$ws.Range("L107").Value = 4741.25
$ws.Range("N107").Value = -8581.25
$ws.Range("H113").Value = 3540
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -8340
$ws.Range("H136").Value = 53427.383
$ws.Range("J136").Value = 41213
$ws.Range("L136").Value = 123639
$ws.Range("N136").Value = -128739

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2904.875
$ws.Range("J7").Value = 3197.75
$ws.Range("L7").Value = 3197.75
$ws.Range("N7").Value = -3421.75
$ws.Range("H9").Value = 4793
$ws.Range("I9").Value = 564.75
$ws.Range("J9").Value = 13249.5
$ws.Range("K9").Value = 564.75
$ws.Range("L9").Value = 13249.5
$ws.Range("M9").Value = -340.75
$ws.Range("N9").Value = -13697.5
$ws.Range("H22").Value = 877.44446
$ws.Range("I22").Value = 727.4
$ws.Range("K22").Value = 727.4
$ws.Range("M22").Value = -432.4
$ws.Range("H25").Value = 16256
$ws.Range("I25").Value = 10004
$ws.Range("J25").Value = 22508
$ws.Range("K25").Value = 10004
$ws.Range("L25").Value = 22508
$ws.Range("M25").Value = -9774
$ws.Range("N25").Value = -22968
$ws.Range("H27").Value = 877.44446
$ws.Range("I27").Value = 727.4
$ws.Range("K27").Value = 727.4
$ws.Range("M27").Value = -620.4
$ws.Range("H46").Value = 2544.9697
$ws.Range("I46").Value = 2388.4062
$ws.Range("K46").Value = 2388.4062
$ws.Range("M46").Value = -2200.4062
$ws.Range("H126").Value = 2904.875
$ws.Range("J126").Value = 3197.75
$ws.Range("L126").Value = 9593.25
$ws.Range("N126").Value = -14533.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17909
$ws.Range("I41").Value = 17909
$ws.Range("K41").Value = 17909
$ws.Range("M41").Value = -17519
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 10000
$ws.Range("K96").Value = 10000
$ws.Range("M96").Value = -8627
$ws.Range("H122").Value = 1455.5
$ws.Range("I122").Value = 1455.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4366.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1916.5
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 23967.334
$ws.Range("J124").Value = 23967.334
$ws.Range("L124").Value = 23967.334
$ws.Range("N124").Value = -33787.334
$ws.Range("H131").Value = 29999.5
$ws.Range("J131").Value = 29999.5
$ws.Range("L131").Value = 29999.5
$ws.Range("N131").Value = -40079.5
$ws.Range("H136").Value = 12529.077
$ws.Range("I136").Value = 11119.208
$ws.Range("K136").Value = 33357.624
$ws.Range("M136").Value = -30807.624
